$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Quantite (column C) updates
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 10
$ws.Range("C5").Value = 1000
$ws.Range("C6").Value = 25
$ws.Range("C11").Value = 31
$ws.Range("C16").Value = 70

# Reference (column B) updates
$ws.Range("B14").Value = 2938475611
$ws.Range("B15").Value = 6574839202
$ws.Range("B16").Value = 3847562911
